$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 14.734287
$ws.Range("H2").Value = 44.202861
$ws.Range("I2").Value = 0.4000023944294819
$ws.Range("J2").Value = 0.400002394429482
$ws.Range("M2").Value = 49.89274333333334
$ws.Range("N2").Value = 149.67823
$ws.Range("O2").Value = 0.8663408689480834
$ws.Range("P2").Value = 0.8663408689480835
$ws.Range("Q2").Value = 735.13399949067
$ws.Range("R2").Value = 6616.20599541603
$ws.Range("S2").Value = 0.3465384219713514
$ws.Range("T2").Value = 0.3465384219713515

$ws.Range("G3").Value = 14.734287
$ws.Range("H3").Value = 44.202861
$ws.Range("I3").Value = 0.4000023944294819
$ws.Range("J3").Value = 0.400002394429482
$ws.Range("O3").Value = 0.06984725491313053
$ws.Range("P3").Value = 0.06984725491313053
$ws.Range("Q3").Value = 59.26892485180799
$ws.Range("R3").Value = 533.420323666272
$ws.Range("S3").Value = 0.02793906920957861
$ws.Range("T3").Value = 0.02793906920957862

$ws.Range("G4").Value = 14.734287
$ws.Range("H4").Value = 44.202861
$ws.Range("I4").Value = 0.4000023944294819
$ws.Range("J4").Value = 0.400002394429482
$ws.Range("M4").Value = 1.266267666666667
$ws.Range("N4").Value = 3.798803
$ws.Range("O4").Value = 0.02198755485004457
$ws.Range("P4").Value = 0.02198755485004457
$ws.Range("Q4").Value = 18.657551219487
$ws.Range("R4").Value = 167.917960975383
$ws.Range("S4").Value = 0.008795074587667396
$ws.Range("T4").Value = 0.008795074587667399

$ws.Range("G5").Value = 14.734287
$ws.Range("H5").Value = 44.202861
$ws.Range("I5").Value = 0.4000023944294819
$ws.Range("J5").Value = 0.400002394429482
$ws.Range("M5").Value = 0.2206823333333333
$ws.Range("N5").Value = 0.6620469999999999
$ws.Range("O5").Value = 0.003831942516052412
$ws.Range("P5").Value = 0.003831942516052413
$ws.Range("Q5").Value = 3.251596835163
$ws.Range("R5").Value = 29.264371516467
$ws.Range("S5").Value = 0.001532786181737098
$ws.Range("T5").Value = 0.001532786181737099

$ws.Range("G6").Value = 14.734287
$ws.Range("H6").Value = 44.202861
$ws.Range("I6").Value = 0.4000023944294819
$ws.Range("J6").Value = 0.400002394429482
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.187988666666667
$ws.Range("N6").Value = 6.563966
$ws.Range("O6").Value = 0.03799237877268909
$ws.Range("P6").Value = 0.03799237877268909
$ws.Range("Q6").Value = 32.238452967414
$ws.Range("R6").Value = 290.146076706726
$ws.Range("S6").Value = 0.01519704247914746
$ws.Range("T6").Value = 0.01519704247914746

$ws.Range("I7").Value = 0.3923645715978801
$ws.Range("J7").Value = 0.3923645715978802
$ws.Range("M7").Value = 49.89274333333334
$ws.Range("N7").Value = 149.67823
$ws.Range("O7").Value = 0.8663408689480834
$ws.Range("P7").Value = 0.8663408689480835
$ws.Range("Q7").Value = 721.0970254030401
$ws.Range("R7").Value = 6489.873228627361
$ws.Range("S7").Value = 0.33992146390255
$ws.Range("T7").Value = 0.3399214639025501

$ws.Range("I8").Value = 0.3923645715978801
$ws.Range("J8").Value = 0.3923645715978802
$ws.Range("O8").Value = 0.06984725491313053
$ws.Range("P8").Value = 0.06984725491313053
$ws.Range("S8").Value = 0.02740558825127839
$ws.Range("T8").Value = 0.0274055882512784

$ws.Range("I9").Value = 0.3923645715978801
$ws.Range("J9").Value = 0.3923645715978802
$ws.Range("M9").Value = 1.266267666666667
$ws.Range("N9").Value = 3.798803
$ws.Range("O9").Value = 0.02198755485004457
$ws.Range("P9").Value = 0.02198755485004457
$ws.Range("Q9").Value = 18.301295675344
$ws.Range("R9").Value = 164.711661078096
$ws.Range("S9").Value = 0.008627137539222627
$ws.Range("T9").Value = 0.008627137539222632

$ws.Range("I10").Value = 0.3923645715978801
$ws.Range("J10").Value = 0.3923645715978802
$ws.Range("M10").Value = 0.2206823333333333
$ws.Range("N10").Value = 0.6620469999999999
$ws.Range("O10").Value = 0.003831942516052412
$ws.Range("P10").Value = 0.003831942516052413
$ws.Range("Q10").Value = 3.189509405456
$ws.Range("R10").Value = 28.705584649104
$ws.Range("S10").Value = 0.001503518483698608
$ws.Range("T10").Value = 0.001503518483698608

$ws.Range("I11").Value = 0.3923645715978801
$ws.Range("J11").Value = 0.3923645715978802
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 2.187988666666667
$ws.Range("N11").Value = 6.563966
$ws.Range("O11").Value = 0.03799237877268909
$ws.Range("P11").Value = 0.03799237877268909
$ws.Range("Q11").Value = 31.622877671968
$ws.Range("R11").Value = 284.605899047712
$ws.Range("S11").Value = 0.01490686342113055
$ws.Range("T11").Value = 0.01490686342113055

$ws.Range("G12").Value = 1.259379333333333
$ws.Range("H12").Value = 3.778138
$ws.Range("I12").Value = 0.03418928576783783
$ws.Range("J12").Value = 0.03418928576783784
$ws.Range("M12").Value = 49.89274333333334
$ws.Range("N12").Value = 149.67823
$ws.Range("O12").Value = 0.8663408689480834
$ws.Range("P12").Value = 0.8663408689480835
$ws.Range("Q12").Value = 62.83388983730445
$ws.Range("R12").Value = 565.5050085357401
$ws.Range("S12").Value = 0.02961957554082297
$ws.Range("T12").Value = 0.02961957554082298

$ws.Range("G13").Value = 1.259379333333333
$ws.Range("H13").Value = 3.778138
$ws.Range("I13").Value = 0.03418928576783783
$ws.Range("J13").Value = 0.03418928576783784
$ws.Range("O13").Value = 0.06984725491313053
$ws.Range("P13").Value = 0.06984725491313053
$ws.Range("Q13").Value = 5.065875197575111
$ws.Range("R13").Value = 45.592876778176
$ws.Range("S13").Value = 0.002388027758324035
$ws.Range("T13").Value = 0.002388027758324035

$ws.Range("G14").Value = 1.259379333333333
$ws.Range("H14").Value = 3.778138
$ws.Range("I14").Value = 0.03418928576783783
$ws.Range("J14").Value = 0.03418928576783784
$ws.Range("M14").Value = 1.266267666666667
$ws.Range("N14").Value = 3.798803
$ws.Range("O14").Value = 0.02198755485004457
$ws.Range("P14").Value = 0.02198755485004457
$ws.Range("Q14").Value = 1.594711329868223
$ws.Range("R14").Value = 14.352401968814
$ws.Range("S14").Value = 0.0007517387961041824
$ws.Range("T14").Value = 0.0007517387961041826

$ws.Range("G15").Value = 1.259379333333333
$ws.Range("H15").Value = 3.778138
$ws.Range("I15").Value = 0.03418928576783783
$ws.Range("J15").Value = 0.03418928576783784
$ws.Range("M15").Value = 0.2206823333333333
$ws.Range("N15").Value = 0.6620469999999999
$ws.Range("O15").Value = 0.003831942516052412
$ws.Range("P15").Value = 0.003831942516052413
$ws.Range("Q15").Value = 0.2779227698317778
$ws.Range("R15").Value = 2.501304928486
$ws.Range("S15").Value = 0.0001310113777272434
$ws.Range("T15").Value = 0.0001310113777272435

$ws.Range("G16").Value = 1.259379333333333
$ws.Range("H16").Value = 3.778138
$ws.Range("I16").Value = 0.03418928576783783
$ws.Range("J16").Value = 0.03418928576783784
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 2.187988666666667
$ws.Range("N16").Value = 6.563966
$ws.Range("O16").Value = 0.03799237877268909
$ws.Range("P16").Value = 0.03799237877268909
$ws.Range("Q16").Value = 2.755507708367555
$ws.Range("R16").Value = 24.799569375308
$ws.Range("S16").Value = 0.001298932294859403
$ws.Range("T16").Value = 0.001298932294859403

$ws.Range("G17").Value = 4.524255666666667
$ws.Range("H17").Value = 13.572767
$ws.Range("I17").Value = 0.1228232556945456
$ws.Range("J17").Value = 0.1228232556945456
$ws.Range("M17").Value = 49.89274333333334
$ws.Range("N17").Value = 149.67823
$ws.Range("O17").Value = 0.8663408689480834
$ws.Range("P17").Value = 0.8663408689480835
$ws.Range("Q17").Value = 225.7275267513789
$ws.Range("R17").Value = 2031.54774076241
$ws.Range("S17").Value = 0.1064068060654453
$ws.Range("T17").Value = 0.1064068060654453

$ws.Range("G18").Value = 4.524255666666667
$ws.Range("H18").Value = 13.572767
$ws.Range("I18").Value = 0.1228232556945456
$ws.Range("J18").Value = 0.1228232556945456
$ws.Range("O18").Value = 0.06984725491313053
$ws.Range("P18").Value = 0.06984725491313053
$ws.Range("Q18").Value = 18.19889683959822
$ws.Range("R18").Value = 163.790071556384
$ws.Range("S18").Value = 0.008578867249757536
$ws.Range("T18").Value = 0.008578867249757538

$ws.Range("G19").Value = 4.524255666666667
$ws.Range("H19").Value = 13.572767
$ws.Range("I19").Value = 0.1228232556945456
$ws.Range("J19").Value = 0.1228232556945456
$ws.Range("M19").Value = 1.266267666666667
$ws.Range("N19").Value = 3.798803
$ws.Range("O19").Value = 0.02198755485004457
$ws.Range("P19").Value = 0.02198755485004457
$ws.Range("Q19").Value = 5.728918666433445
$ws.Range("R19").Value = 51.56026799790101
$ws.Range("S19").Value = 0.002700583071444869
$ws.Range("T19").Value = 0.00270058307144487

$ws.Range("G20").Value = 4.524255666666667
$ws.Range("H20").Value = 13.572767
$ws.Range("I20").Value = 0.1228232556945456
$ws.Range("J20").Value = 0.1228232556945456
$ws.Range("M20").Value = 0.2206823333333333
$ws.Range("N20").Value = 0.6620469999999999
$ws.Range("O20").Value = 0.003831942516052412
$ws.Range("P20").Value = 0.003831942516052413
$ws.Range("Q20").Value = 0.9984232971165555
$ws.Range("R20").Value = 8.985809674048999
$ws.Range("S20").Value = 0.0004706516554559057
$ws.Range("T20").Value = 0.0004706516554559059

$ws.Range("G21").Value = 4.524255666666667
$ws.Range("H21").Value = 13.572767
$ws.Range("I21").Value = 0.1228232556945456
$ws.Range("J21").Value = 0.1228232556945456
$ws.Range("K21").Value = 3
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 2.187988666666667
$ws.Range("N21").Value = 6.563966
$ws.Range("O21").Value = 0.03799237877268909
$ws.Range("P21").Value = 0.03799237877268909
$ws.Range("Q21").Value = 9.899020123769111
$ws.Range("R21").Value = 89.091181113922
$ws.Range("S21").Value = 0.004666347652442017
$ws.Range("T21").Value = 0.004666347652442018

$ws.Range("G22").Value = 1.864631
$ws.Range("H22").Value = 5.593893
$ws.Range("I22").Value = 0.05062049251025444
$ws.Range("J22").Value = 0.05062049251025445
$ws.Range("M22").Value = 49.89274333333334
$ws.Range("N22").Value = 149.67823
$ws.Range("O22").Value = 0.8663408689480834
$ws.Range("P22").Value = 0.8663408689480835
$ws.Range("Q22").Value = 93.03155589437667
$ws.Range("R22").Value = 837.28400304939
$ws.Range("S22").Value = 0.04385460146791378
$ws.Range("T22").Value = 0.04385460146791379

$ws.Range("G23").Value = 1.864631
$ws.Range("H23").Value = 5.593893
$ws.Range("I23").Value = 0.05062049251025444
$ws.Range("J23").Value = 0.05062049251025445
$ws.Range("O23").Value = 0.06984725491313053
$ws.Range("P23").Value = 0.06984725491313053
$ws.Range("Q23").Value = 7.500510517770666
$ws.Range("R23").Value = 67.50459465993599
$ws.Range("S23").Value = 0.003535702444191957
$ws.Range("T23").Value = 0.003535702444191957

$ws.Range("G24").Value = 1.864631
$ws.Range("H24").Value = 5.593893
$ws.Range("I24").Value = 0.05062049251025444
$ws.Range("J24").Value = 0.05062049251025445
$ws.Range("M24").Value = 1.266267666666667
$ws.Range("N24").Value = 3.798803
$ws.Range("O24").Value = 0.02198755485004457
$ws.Range("P24").Value = 0.02198755485004457
$ws.Range("Q24").Value = 2.361121945564333
$ws.Range("R24").Value = 21.250097510079
$ws.Range("S24").Value = 0.00111302085560549
$ws.Range("T24").Value = 0.00111302085560549

$ws.Range("G25").Value = 1.864631
$ws.Range("H25").Value = 5.593893
$ws.Range("I25").Value = 0.05062049251025444
$ws.Range("J25").Value = 0.05062049251025445
$ws.Range("M25").Value = 0.2206823333333333
$ws.Range("N25").Value = 0.6620469999999999
$ws.Range("O25").Value = 0.003831942516052412
$ws.Range("P25").Value = 0.003831942516052413
$ws.Range("Q25").Value = 0.4114911198856666
$ws.Range("R25").Value = 3.703420078970999
$ws.Range("S25").Value = 0.0001939748174335567
$ws.Range("T25").Value = 0.0001939748174335567

$ws.Range("G26").Value = 1.864631
$ws.Range("H26").Value = 5.593893
$ws.Range("I26").Value = 0.05062049251025444
$ws.Range("J26").Value = 0.05062049251025445
$ws.Range("K26").Value = 3
$ws.Range("L26").Value = 1
$ws.Range("M26").Value = 2.187988666666667
$ws.Range("N26").Value = 6.563966
$ws.Range("O26").Value = 0.03799237877268909
$ws.Range("P26").Value = 0.03799237877268909
$ws.Range("Q26").Value = 4.079791495515333
$ws.Range("R26").Value = 36.71812345963799
$ws.Range("S26").Value = 0.001923192925109658
$ws.Range("T26").Value = 0.001923192925109658
